$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45192 = 2023-09-23) for every
# data row (rows 2 through 472). Update it to serial 45202 (2023-10-03), i.e.
# shift it 10 days forward, for all data rows.
$lastRow = 472
$ws.Range("C2:C$lastRow").Value = 45202
